$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before G (shifts old G..  to H)
$ws.Columns("G:G").Insert()
$ws.Range("G2:G251").Clear()

# Set header for new column G1
$ws.Cells.Item(1,7).Value = "Distribuidor"
